# Fix semantic nonsense in questions: Expand KB and fix fallbacks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Row 2 ---
$ws.Range("F2").Value = "Kendrick Lamar"
$ws.Range("K2").Value = "101,50,54"

# --- Row 3 ---
$ws.Range("D3").Value = "Who performed at the 2021 Super Bowl Halftime Show?"
$ws.Range("G3").Value = "Justin Bieber"
$ws.Range("K3").Value = "10,100,50,54"

# --- Row 4 ---
$ws.Range("F4").Value = "Maluma"
$ws.Range("G4").Value = "Daddy Yankee"
$ws.Range("K4").Value = "50,54,99"

# --- Row 5 ---
$ws.Range("D5").Value = "Who sings 'thank u, next'?"
$ws.Range("F5").Value = "Taylor Swift"
$ws.Range("G5").Value = "Selena Gomez"
$ws.Range("K5").Value = "10,100,101,50,99"

# --- Row 6 ---
$ws.Range("F6").Value = "Titanic"
$ws.Range("G6").Value = "Star Wars"
$ws.Range("K6").Value = "10,101,51"

# --- Row 7 ---
$ws.Range("D7").Value = "In which movie does Iron Man snap his fingers?"
$ws.Range("G7").Value = "Age of Ultron"
$ws.Range("K7").Value = "51,54,99"

# --- Row 8 ---
$ws.Range("F8").Value = "Amazon"
$ws.Range("G8").Value = "Twitter"
$ws.Range("K8").Value = "100,101,54,7,99"

# --- Row 9 ---
$ws.Range("K9").Value = "10,101,99"

# --- Row 10 ---
$ws.Range("K10").Value = "10,100,101,54,99"

# --- Add new "tags" sheet after "data" ---
$tags = $wb.Worksheets.Add($null, $ws)
$tags.Name = "tags"

$tags.Range("A1").Value = "TAG_ID"
$tags.Range("B1").Value = "TAG_NAME"
$tags.Range("C1").Value = "CATEGORY"

# Match the bold/bordered header style used on the "data" sheet's row 1.
$ws.Range("A1:C1").Copy()
$tags.Range("A1:C1").PasteSpecial(-4122)

# TAG_ID values are text-typed (e.g. "4", "50") in the source data, so
# force the Text number format before assignment to avoid numeric coercion.
$tagRows = @(
    @("4",   ":COUNTRY",      "Geography"),
    @("7",   ":COMPANY",      "Business"),
    @("10",  "ED:HISTORY",    "Education"),
    @("50",  "EN:MUSIC",      "Entertainment"),
    @("51",  "EN:MOVIE",      "Entertainment"),
    @("54",  "EN:Facts",      "General"),
    @("99",  "NEW:Viral",     "Trends"),
    @("100", "NEW:Tech",      "Technology"),
    @("101", "NEW:Nostalgia", "Lifestyle")
)

$r = 2
foreach ($row in $tagRows) {
    $idCell = $tags.Cells.Item($r, 1)
    $idCell.NumberFormat = "@"
    $idCell.Value = $row[0]
    $tags.Cells.Item($r, 2).Value = $row[1]
    $tags.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Activate()
